$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1153286666666667
$ws.Range("H2").Value = 0.345986
$ws.Range("I2").Value = 0.1212009326543938
$ws.Range("J2").Value = 0.1212009326543938
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1104363333333333
$ws.Range("N2").Value = 0.331309
$ws.Range("O2").Value = 0.2546765664720067
$ws.Range("P2").Value = 0.2546765664720067
$ws.Range("Q2").Value = 0.01273647507488889
$ws.Range("R2").Value = 0.114628275674
$ws.Range("S2").Value = 0.03086703738162594
$ws.Range("T2").Value = 0.03086703738162594
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1153286666666667
$ws.Range("H3").Value = 0.345986
$ws.Range("I3").Value = 0.1212009326543938
$ws.Range("J3").Value = 0.1212009326543938
$ws.Range("O3").Value = 0.6439575340475562
$ws.Range("P3").Value = 0.6439575340475563
$ws.Range("Q3").Value = 0.03220456909444445
$ws.Range("R3").Value = 0.28984112185
$ws.Range("S3").Value = 0.07804825371638738
$ws.Range("T3").Value = 0.07804825371638738
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1153286666666667
$ws.Range("H4").Value = 0.345986
$ws.Range("I4").Value = 0.1212009326543938
$ws.Range("J4").Value = 0.1212009326543938
$ws.Range("O4").Value = 0.101365899480437
$ws.Range("P4").Value = 0.101365899480437
$ws.Range("Q4").Value = 0.005069348429111112
$ws.Range("R4").Value = 0.04562413586200001
$ws.Range("S4").Value = 0.0122856415563805
$ws.Range("T4").Value = 0.0122856415563805
$ws.Range("I5").Value = 0.5471739422864045
$ws.Range("J5").Value = 0.5471739422864045
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1104363333333333
$ws.Range("N5").Value = 0.331309
$ws.Range("O5").Value = 0.2546765664720067
$ws.Range("P5").Value = 0.2546765664720067
$ws.Range("Q5").Value = 0.05750011262233334
$ws.Range("R5").Value = 0.5175010136010001
$ws.Range("S5").Value = 0.1393523808844535
$ws.Range("T5").Value = 0.1393523808844535
$ws.Range("I6").Value = 0.5471739422864045
$ws.Range("J6").Value = 0.5471739422864045
$ws.Range("O6").Value = 0.6439575340475562
$ws.Range("P6").Value = 0.6439575340475563
$ws.Range("S6").Value = 0.3523567825698328
$ws.Range("T6").Value = 0.3523567825698329
$ws.Range("I7").Value = 0.5471739422864045
$ws.Range("J7").Value = 0.5471739422864045
$ws.Range("O7").Value = 0.101365899480437
$ws.Range("P7").Value = 0.101365899480437
$ws.Range("S7").Value = 0.05546477883211812
$ws.Range("T7").Value = 0.05546477883211812
$ws.Range("H8").Value = 0.9466730000000001
$ws.Range("I8").Value = 0.3316251250592017
$ws.Range("J8").Value = 0.3316251250592017
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1104363333333333
$ws.Range("N8").Value = 0.331309
$ws.Range("O8").Value = 0.2546765664720067
$ws.Range("P8").Value = 0.2546765664720067
$ws.Range("Q8").Value = 0.03484903166188889
$ws.Range("R8").Value = 0.3136412849570001
$ws.Range("S8").Value = 0.08445714820592733
$ws.Range("T8").Value = 0.08445714820592733
$ws.Range("H9").Value = 0.9466730000000001
$ws.Range("I9").Value = 0.3316251250592017
$ws.Range("J9").Value = 0.3316251250592017
$ws.Range("O9").Value = 0.6439575340475562
$ws.Range("P9").Value = 0.6439575340475563
$ws.Range("S9").Value = 0.213552497761336
$ws.Range("T9").Value = 0.213552497761336
$ws.Range("H10").Value = 0.9466730000000001
$ws.Range("I10").Value = 0.3316251250592017
$ws.Range("J10").Value = 0.3316251250592017
$ws.Range("O10").Value = 0.101365899480437
$ws.Range("P10").Value = 0.101365899480437
$ws.Range("S10").Value = 0.0336154790919384
$ws.Range("T10").Value = 0.0336154790919384
